$d = $word.ActiveDocument

# 1) Fix joined words: "Metode" + "Pengumpulan Data" -> "Metode Pengumpulan Data"
#    (the Heading2 "Metode Pengumpulan Data" title currently has no space
#    between the two runs, so the rendered text reads "MetodePengumpulan Data")
$d.Content.Find.Execute("MetodePengumpulan Data", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Metode Pengumpulan Data", 2)

# 2) Fix the cached page-number field result in the header from 25 to 22
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    if ($hdr.Exists) {
        $hdr.Range.Find.Execute("25", $true, $false, $false, $false, $false,
                                 $true, 1, $false, "22", 2)
    }
}
